# Apply weekly update to "Fruta, Terminal Hortofrutícola Agro Chillán - Arándano (blue)"
# Rows 2-9 hold the rolling weekly price records; this commit adds a new
# week (row 9 here becomes the new record pushed in) and shifts the other
# weekly snapshots, matching the source feed row order exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(2, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(2, 3).Value = "Ñuble"
$ws.Cells.Item(2, 4).Value = 44187
$ws.Cells.Item(2, 5).Value = 16
$ws.Cells.Item(2, 6).Value = "Fruta"
$ws.Cells.Item(2, 7).Value = 100101
$ws.Cells.Item(2, 8).Value = "Berries"
$ws.Cells.Item(2, 9).Value = 100101001
$ws.Cells.Item(2, 10).Value = "Arándano (blue)"
$ws.Cells.Item(2, 11).Value = "Sin especificar"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 80
$ws.Cells.Item(2, 14).Value = 2800
$ws.Cells.Item(2, 15).Value = 3000
$ws.Cells.Item(2, 16).Value = 2900
$ws.Cells.Item(2, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(2, 18).Value = "Provincia de Linares"
$ws.Cells.Item(2, 19).Value = 1450
$ws.Cells.Item(2, 20).Value = 2

# Row 3
$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value = "Ñuble"
$ws.Cells.Item(3, 4).Value = 44187
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100101
$ws.Cells.Item(3, 8).Value = "Berries"
$ws.Cells.Item(3, 9).Value = 100101001
$ws.Cells.Item(3, 10).Value = "Arándano (blue)"
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 65
$ws.Cells.Item(3, 14).Value = 1400
$ws.Cells.Item(3, 15).Value = 1500
$ws.Cells.Item(3, 16).Value = 1446
$ws.Cells.Item(3, 17).Value = "`$/envase 1 kilo"
$ws.Cells.Item(3, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(3, 19).Value = 1446
$ws.Cells.Item(3, 20).Value = 1

# Row 4
$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(4, 3).Value = "Ñuble"
$ws.Cells.Item(4, 4).Value = 44181
$ws.Cells.Item(4, 5).Value = 16
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100101
$ws.Cells.Item(4, 8).Value = "Berries"
$ws.Cells.Item(4, 9).Value = 100101001
$ws.Cells.Item(4, 10).Value = "Arándano (blue)"
$ws.Cells.Item(4, 11).Value = "Sin especificar"
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 65
$ws.Cells.Item(4, 14).Value = 3600
$ws.Cells.Item(4, 15).Value = 3800
$ws.Cells.Item(4, 16).Value = 3692
$ws.Cells.Item(4, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(4, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(4, 19).Value = 1846
$ws.Cells.Item(4, 20).Value = 2

# Row 5
$ws.Cells.Item(5, 1).Value = 7
$ws.Cells.Item(5, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(5, 3).Value = "Ñuble"
$ws.Cells.Item(5, 4).Value = 44181
$ws.Cells.Item(5, 5).Value = 16
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100101
$ws.Cells.Item(5, 8).Value = "Berries"
$ws.Cells.Item(5, 9).Value = 100101001
$ws.Cells.Item(5, 10).Value = "Arándano (blue)"
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 80
$ws.Cells.Item(5, 14).Value = 1800
$ws.Cells.Item(5, 15).Value = 2000
$ws.Cells.Item(5, 16).Value = 1875
$ws.Cells.Item(5, 17).Value = "`$/envase 1 kilo"
$ws.Cells.Item(5, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(5, 19).Value = 1875
$ws.Cells.Item(5, 20).Value = 1

# Row 6
$ws.Cells.Item(6, 1).Value = 7
$ws.Cells.Item(6, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(6, 3).Value = "Ñuble"
$ws.Cells.Item(6, 4).Value = 44594
$ws.Cells.Item(6, 5).Value = 16
$ws.Cells.Item(6, 6).Value = "Fruta"
$ws.Cells.Item(6, 7).Value = 100101
$ws.Cells.Item(6, 8).Value = "Berries"
$ws.Cells.Item(6, 9).Value = 100101001
$ws.Cells.Item(6, 10).Value = "Arándano (blue)"
$ws.Cells.Item(6, 11).Value = "Sin especificar"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 120
$ws.Cells.Item(6, 14).Value = 2500
$ws.Cells.Item(6, 15).Value = 2800
$ws.Cells.Item(6, 16).Value = 2650
$ws.Cells.Item(6, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(6, 18).Value = "Provincia de Linares"
$ws.Cells.Item(6, 19).Value = 1325
$ws.Cells.Item(6, 20).Value = 2

# Row 7
$ws.Cells.Item(7, 1).Value = 7
$ws.Cells.Item(7, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(7, 3).Value = "Ñuble"
$ws.Cells.Item(7, 4).Value = 44540
$ws.Cells.Item(7, 5).Value = 16
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100101
$ws.Cells.Item(7, 8).Value = "Berries"
$ws.Cells.Item(7, 9).Value = 100101001
$ws.Cells.Item(7, 10).Value = "Arándano (blue)"
$ws.Cells.Item(7, 11).Value = "Sin especificar"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 240
$ws.Cells.Item(7, 14).Value = 3500
$ws.Cells.Item(7, 15).Value = 3800
$ws.Cells.Item(7, 16).Value = 3650
$ws.Cells.Item(7, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(7, 18).Value = "Región del Maule"
$ws.Cells.Item(7, 19).Value = 1825
$ws.Cells.Item(7, 20).Value = 2

# Row 8
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8, 3).Value = "Ñuble"
$ws.Cells.Item(8, 4).Value = 44539
$ws.Cells.Item(8, 5).Value = 16
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100101
$ws.Cells.Item(8, 8).Value = "Berries"
$ws.Cells.Item(8, 9).Value = 100101001
$ws.Cells.Item(8, 10).Value = "Arándano (blue)"
$ws.Cells.Item(8, 11).Value = "Sin especificar"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 200
$ws.Cells.Item(8, 14).Value = 3800
$ws.Cells.Item(8, 15).Value = 4000
$ws.Cells.Item(8, 16).Value = 3900
$ws.Cells.Item(8, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(8, 18).Value = "Región del Maule"
$ws.Cells.Item(8, 19).Value = 1950
$ws.Cells.Item(8, 20).Value = 2

# Row 9
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44174
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100101
$ws.Cells.Item(9, 8).Value = "Berries"
$ws.Cells.Item(9, 9).Value = 100101001
$ws.Cells.Item(9, 10).Value = "Arándano (blue)"
$ws.Cells.Item(9, 11).Value = "Sin especificar"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 150
$ws.Cells.Item(9, 14).Value = 3700
$ws.Cells.Item(9, 15).Value = 3800
$ws.Cells.Item(9, 16).Value = 3747
$ws.Cells.Item(9, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(9, 18).Value = "Provincia de Linares"
$ws.Cells.Item(9, 19).Value = 1874
$ws.Cells.Item(9, 20).Value = 2

# Ensure the date column keeps the workbook-standard datetime format
# (rows 2-8 already carry it; row 9 is newly created by this edit).
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
